# Update BIBI retention metrics data (atualizei dados da BIBI)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 989
$ws.Range("D2").Value = 989
$ws.Range("C3").Value = 106
$ws.Range("D3").Value = 989
$ws.Range("E3").Value = 0.1071789686552073
$ws.Range("C4").Value = 37
$ws.Range("D4").Value = 989
$ws.Range("E4").Value = 0.03741152679474216
$ws.Range("D5").Value = 989
$ws.Range("E5").Value = 0.005055611729019211
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 989
$ws.Range("E6").Value = 0.006066734074823054
$ws.Range("D7").Value = 989
$ws.Range("E7").Value = 0.007077856420626896
$ws.Range("D8").Value = 989
$ws.Range("E8").Value = 0.004044489383215369
$ws.Range("D9").Value = 989
$ws.Range("E9").Value = 0.002022244691607685
$ws.Range("C10").Value = 2107
$ws.Range("D10").Value = 2107
$ws.Range("C11").Value = 129
$ws.Range("D11").Value = 2107
$ws.Range("E11").Value = 0.06122448979591837
$ws.Range("D12").Value = 2107
$ws.Range("E12").Value = 0.01756051257712387
$ws.Range("D13").Value = 2107
$ws.Range("E13").Value = 0.01613668723303275
$ws.Range("D14").Value = 2107
$ws.Range("E14").Value = 0.0132890365448505
$ws.Range("D15").Value = 2107
$ws.Range("E15").Value = 0.01044138585666825
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 2107
$ws.Range("E16").Value = 0.003322259136212625
$ws.Range("C17").Value = 2654
$ws.Range("D17").Value = 2654
$ws.Range("C18").Value = 167
$ws.Range("D18").Value = 2654
$ws.Range("E18").Value = 0.06292388847023361
$ws.Range("D19").Value = 2654
$ws.Range("E19").Value = 0.04860587792012058
$ws.Range("D20").Value = 2654
$ws.Range("E20").Value = 0.02901281085154484
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 2654
$ws.Range("E21").Value = 0.02147701582516956
$ws.Range("C22").Value = 22
$ws.Range("D22").Value = 2654
$ws.Range("E22").Value = 0.008289374529012811
$ws.Range("C23").Value = 2252
$ws.Range("D23").Value = 2252
$ws.Range("C24").Value = 268
$ws.Range("D24").Value = 2252
$ws.Range("E24").Value = 0.1190053285968028
$ws.Range("C25").Value = 127
$ws.Range("D25").Value = 2252
$ws.Range("E25").Value = 0.0563943161634103
$ws.Range("C26").Value = 93
$ws.Range("D26").Value = 2252
$ws.Range("E26").Value = 0.04129662522202487
$ws.Range("C27").Value = 34
$ws.Range("D27").Value = 2252
$ws.Range("E27").Value = 0.01509769094138544
$ws.Range("C28").Value = 2312
$ws.Range("D28").Value = 2312
$ws.Range("C29").Value = 239
$ws.Range("D29").Value = 2312
$ws.Range("E29").Value = 0.1033737024221453
$ws.Range("C30").Value = 131
$ws.Range("D30").Value = 2312
$ws.Range("E30").Value = 0.05666089965397924
$ws.Range("C31").Value = 34
$ws.Range("D31").Value = 2312
$ws.Range("E31").Value = 0.01470588235294118
$ws.Range("C33").Value = 201
$ws.Range("E33").Value = 0.08909574468085106
$ws.Range("C34").Value = 53
$ws.Range("E34").Value = 0.02349290780141844
$ws.Range("C35").Value = 1930
$ws.Range("D35").Value = 1930
$ws.Range("C36").Value = 95
$ws.Range("D36").Value = 1930
$ws.Range("E36").Value = 0.04922279792746114
$ws.Range("C37").Value = 514
$ws.Range("D37").Value = 514

Write-Output "Updated BIBI retention metrics cells."
